# Updated symbol list on Sun Jan 29 21:48:18 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'317.86"
$ws.Range("E2").Value = "'3.93%"
$ws.Range("D3").Value = "'39.73"
$ws.Range("E3").Value = "'2.22%"
$ws.Range("D4").Value = "'5.134"
$ws.Range("E4").Value = "'0.49%"
$ws.Range("D5").Value = "'0.08209"
$ws.Range("E5").Value = "'1.59%"
$ws.Range("D6").Value = "'2.062"
$ws.Range("E6").Value = "'6.73%"
$ws.Range("D7").Value = "'8.354"
$ws.Range("D8").Value = "'4.304"
$ws.Range("E8").Value = "'2.44%"
$ws.Range("D9").Value = "'0.9390"
$ws.Range("E9").Value = "'1.27%"
$ws.Range("D10").Value = "'0.1359"
$ws.Range("E10").Value = "'-6.46%"
$ws.Range("D11").Value = "'0.1987"
$ws.Range("E11").Value = "'3.92%"
$ws.Range("D12").Value = "'0.09169"
$ws.Range("E12").Value = "'1.10%"
$ws.Range("D13").Value = "'0.03504"
$ws.Range("E13").Value = "'-0.06%"
$ws.Range("D14").Value = "'0.09786"
$ws.Range("E14").Value = "'0.12%"
$ws.Range("E15").Value = "'1.04%"
$ws.Range("D16").Value = "'0.006195"
$ws.Range("E16").Value = "'4.64%"
$ws.Range("D17").Value = "'3.681"
$ws.Range("E17").Value = "'-2.45%"
$ws.Range("D18").Value = "'3.237"
$ws.Range("E18").Value = "'-5.04%"
$ws.Range("D19").Value = "'0.3482"
$ws.Range("E19").Value = "'0.56%"
$ws.Range("E20").Value = "'-0.53%"
$ws.Range("D21").Value = "'4.990"
$ws.Range("E21").Value = "'6.22%"
$ws.Range("E22").Value = "'1.29%"
$ws.Range("D23").Value = "'0.04349"
$ws.Range("D24").Value = "'0.001228"
$ws.Range("E24").Value = "'-0.36%"
$ws.Range("D25").Value = "'0.004812"
$ws.Range("E25").Value = "'12.56%"
$ws.Range("E26").Value = "'-0.08%"
$ws.Range("D27").Value = "'0.0003996"
$ws.Range("E27").Value = "'-10.14%"
$ws.Range("D39").Value = "'0.02253"
$ws.Range("E39").Value = "'10.72%"
$ws.Range("D40").Value = "'0.05198"
$ws.Range("E40").Value = "'2.82%"
$ws.Range("D41").Value = "'0.007745"
$ws.Range("E41").Value = "'2.88%"
$ws.Range("D42").Value = "'0.009856"
$ws.Range("E42").Value = "'1.46%"
$ws.Range("D43").Value = "'0.1407"
$ws.Range("E43").Value = "'4.81%"
$ws.Range("E44").Value = "'-2.84%"
$ws.Range("D45").Value = "'0.009671"
$ws.Range("E45").Value = "'-2.47%"
$ws.Range("D46").Value = "'0.00006604"
$ws.Range("E46").Value = "'6.54%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.30%"
$ws.Range("D48").Value = "'0.001688"
$ws.Range("E48").Value = "'-6.40%"
$ws.Range("D49").Value = "'0.002944"
$ws.Range("E49").Value = "'2.45%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.30%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.30%"
